$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.548.77"
$ws.Range("E2").Value = "  -0.41%  "
$ws.Range("D3").Value = "2.081.45"
$ws.Range("E3").Value = "  +0.31%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'233.25"
$ws.Range("E5").Value = "  -0.11%  "
$ws.Range("D6").Value = "'0.635"
$ws.Range("E6").Value = "  +2.03%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").Value = "'57.88"
$ws.Range("E8").Value = "  -0.31%  "
$ws.Range("D9").Value = "'0.391"
$ws.Range("E9").Value = "  -0.90%  "
$ws.Range("D10").Value = "'0.0779"
$ws.Range("E10").Value = "  -0.52%  "
$ws.Range("E11").Value = "  +2.88%  "
$ws.Range("D12").Value = "'15.08"
$ws.Range("E12").Value = "  +2.24%  "
$ws.Range("D13").Value = "2.388.64"
$ws.Range("E13").Value = "  +0.26%  "
$ws.Range("D14").Value = "'21.06"
$ws.Range("E14").Value = "  +1.18%  "
$ws.Range("D15").Value = "'0.772"
$ws.Range("E15").Value = "  -0.09%  "
$ws.Range("D16").Value = "'5.34"
$ws.Range("E16").Value = "  +0.68%  "
$ws.Range("D17").Value = "2.091.85"
$ws.Range("E17").Value = "  -3.16%  "
$ws.Range("D18").Value = "37.510.90"
$ws.Range("E18").Value = "  -0.36%  "
$ws.Range("B19").Value = "Litecoin"
$ws.Range("C19").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D19").Value = "'70.76"
$ws.Range("E19").Value = "  -0.36%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "'6.03"
$ws.Range("E20").Value = "  -2.20%  "
$ws.Range("D21").Value = "0.0₃0832"
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("D22").Value = "'228.98"
$ws.Range("E22").Value = "  +0.38%  "
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("E24").Value = "  -0.94%  "
$ws.Range("D25").Value = "'2.38"
$ws.Range("E25").Value = "  -0.34%  "
$ws.Range("D26").Value = "'9.67"
$ws.Range("E26").Value = "  +7.10%  "
$ws.Range("D27").Value = "'170.54"
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("E28").Value = "  -3.83%  "
$ws.Range("D29").Value = "'19.50"
$ws.Range("E29").Value = "  +0.40%  "
$ws.Range("E30").Value = "  +0.30%  "
$ws.Range("E31").Value = "  +0.89%  "
$ws.Range("D32").Value = "'4.65"
$ws.Range("E32").Value = "  -0.40%  "
$ws.Range("E33").Value = "  +1.58%  "
$ws.Range("D34").Value = "'4.65"
$ws.Range("E34").Value = "  +0.47%  "
$ws.Range("D35").Value = "'2.48"
$ws.Range("E35").Value = "  +0.20%  "
$ws.Range("E36").Value = "  -0.69%  "
$ws.Range("D37").Value = "'3.32"
$ws.Range("E37").Value = "  -1.86%  "
$ws.Range("E38").Value = "  +0.05%  "
$ws.Range("D39").Value = "'5.33"
$ws.Range("E39").Value = "  +0.89%  "
$ws.Range("E40").Value = "  +8.14%  "
$ws.Range("D41").Value = "'100.53"
$ws.Range("E41").Value = "  +2.48%  "
$ws.Range("E42").Value = "  -1.05%  "
$ws.Range("E43").Value = "  +0.70%  "
$ws.Range("D45").Value = "'16.80"
$ws.Range("E45").Value = "  +2.56%  "
$ws.Range("D46").Value = "1.460.31"
$ws.Range("E46").Value = "  +0.91%  "
$ws.Range("E47").Value = "  -0.70%  "
$ws.Range("D48").Value = "'3.99"
$ws.Range("E48").Value = "  -4.77%  "
$ws.Range("D49").Value = "'7.26"
$ws.Range("E49").Value = "  -1.84%  "
$ws.Range("E50").Value = "  -1.90%  "
$ws.Range("D51").Value = "2.273.28"
$ws.Range("E51").Value = "  +0.33%  "
